$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.143.89"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "1.851.40"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6875"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07772"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3040"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("D12").Value = "1.843.20"
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.201"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "29.145.69"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.735"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "234.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "2.101.63"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.962"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1425"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.965"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.405"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.520"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.004"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.030"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7035"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.677"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9166"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.72%  "
$ws.Range("D42").Value = "1.098.54"
$ws.Range("E42").Value = "  +5.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.020"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("D49").Value = "1.997.54"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.140"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.916"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.22%  "
